$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date updated
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now carries a value
$ws.Range("B9").Value = "Alvearie Team"

# The two duplicate "Contact" / "No display for ContactDetail" rows (10 and 11)
# are replaced with a single "Jurisdiction" row, so delete one of them first
# (this shifts every following row up by one).
$ws.Rows.Item(11).Delete()

# Row 10 becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# "Case Sensitive" row now carries the value "true" (kept as text, not boolean)
$ws.Range("B14").Value = "'true"

# Dimension should now read A1:B21 (handled automatically as content now ends at row 21)
